# Update nomenclature. Clean code.
# Applies the content changes described by the commit:
#  - kite sheet: rename "obgen.p" -> "obGen.p", remove the two obsolete
#    "prop.p" rows, and update the active-cell selection
#  - tether sheet: rename "sigma" -> "sigma_max", update the active-cell
#    selection, and make this sheet the active tab
#  - metrics sheet: update the active-cell selection

$wb = $excel.ActiveWorkbook

# --- kite sheet -----------------------------------------------------------
$kite = $wb.Worksheets.Item("kite")
$kite.Range("A7").Value2 = "obGen.p"
$kite.Rows("8:9").Delete()
[void]$kite.Range("B15").Select()

# --- metrics sheet ----------------------------------------------------------
$metrics = $wb.Worksheets.Item("metrics")
[void]$metrics.Range("H11").Select()

# --- tether sheet -----------------------------------------------------------
# Selecting/activating this sheet last makes it the workbook's active tab,
# matching the target (tether becomes tabSelected / workbook activeTab).
$tether = $wb.Worksheets.Item("tether")
$tether.Range("A9").Value2 = "sigma_max"
[void]$tether.Range("E16").Select()
[void]$tether.Activate()
